$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything currently in A:H shifts to B:I
$ws.Columns.Item(1).Insert()

# Header for the new column
$ws.Range("A1").Value = "id"

# Fill sequential ids 1..30 for the data rows (rows 2..31)
for ($i = 2; $i -le 31; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

$ws.Range("E7").Select()
